# Test data for Delegate
# Insert a new "DelegateTo" column (O) on the "Transmittals_New" sheet,
# shifting the existing "Action-Level3" column (formerly O) to P, and add
# a new data row (row 9) that exercises the Delegate scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transmittals_New")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Insert a new column before the current column O (15) - this pushes
# "Action-Level3" (and everything further right) one column to the right,
# and the new column inherits the width of its left neighbour (N), matching
# the target file.
$nWidth = $ws.Columns("N").ColumnWidth
$ws.Columns("O").Insert()
$ws.Columns("O").ColumnWidth = $nWidth

# New header for the inserted column.
$ws.Range("O1").Value = "DelegateTo"

# Existing rows 4-6 already carried values in what is now column P
# (Action-Level3); those values simply shifted right automatically with
# the column insert, so nothing else to touch there.

# Append the new test-data row (row 9) for the Delegate scenario.
$ws.Range("A9").Value = "AutoTestAdmin"
$ws.Range("C9").Value = "New Transmittal from Automation"
$ws.Range("D9").Value = "UnTick"
$ws.Range("E9").Value = "Change Note"
$ws.Range("F9").Value = "Issued for Approval"
$ws.Range("M9").Value = "Delegate"
$ws.Range("L9").Value = "Delegate- Message for New transmittal"
$ws.Range("O9").Value = "AutoTestUser"
$ws.Range("P9").Value = "Approved"

# Restore the view to tabSelected with the scroll/selection state observed
# in the target workbook.
$ws.Activate()
$ws.Range("L12").Select()
